# Capacity Supply Curve.xlsx — "Updated CSC variables for calibration."
#
# 1) Sheet "CSC-CSCCCMvSoECBtY": recalibrate the Cost multiplier curve (row 2,
#    columns C:N) with new calibrated values.
# 2) Sheet "CSC-CSCSoCECBiaSY": recalibrate "Share of existing capacity"
#    table — every technology row moves from 0.33 to 0.50, except the
#    "natural gas combined cycle es" row (row 7) which moves from 0.50 to
#    0.75. Rows that were already 0 (no existing capacity) stay untouched.

$wb = $excel.ActiveWorkbook

# --- Sheet 2: CSC-CSCCCMvSoECBtY (Cost multiplier vs Share of Existing Capacity) ---
$wsCM = $wb.Worksheets.Item("CSC-CSCCCMvSoECBtY")

$wsCM.Range("C2").Value = 1.1427461300794932
$wsCM.Range("D2").Value = 1.4438453988846189
$wsCM.Range("E2").Value = 1.9309364276944523
$wsCM.Range("F2").Value = 2.6390158215457884
$wsCM.Range("G2").Value = 3.6067497647680336
$wsCM.Range("H2").Value = 4.8742339619263264
$wsCM.Range("I2").Value = 6.4818210260626374
$wsCM.Range("J2").Value = 8.4695622497729151
$wsCM.Range("K2").Value = 10.876965857390774
$wsCM.Range("L2").Value = 13.742909902417514
$wsCM.Range("M2").Value = 17.105628240207455
$wsCM.Range("N2").Value = 21.002729596824665

$wsCM.Range("C3:N3").Select() | Out-Null

# --- Sheet 4: CSC-CSCSoCECBiaSY (Share of Cost Effective Capacity Built in a Single Year) ---
$wsSC = $wb.Worksheets.Item("CSC-CSCSoCECBiaSY")

$rowsToHalf = @(2,3,4,5,6,8,9,10,11,12,13,14,15,18,19,20,21,22,23,24,25)
foreach ($r in $rowsToHalf) {
    $wsSC.Range("B$r`:AE$r").Value = 0.5
}

$wsSC.Range("B7:AE7").Value = 0.75

$wsSC.Activate() | Out-Null
$wsSC.Range("B7:AE7").Select() | Out-Null
